$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows (Customer, Project, Package, Commit number, Remarks) ---
# Row 7 - qtbrowser-2.0.9
$ws.Range("A7").Value = "UPC"
$ws.Range("B7").Value = "DAWN Europe- Browser"
$ws.Range("C7").Value = "qtbrowser-2.0.9"
$ws.Range("D7").Value = "a3e844872639992072f8f81830f67e8a56a592f0"
$ws.Range("E7").Value = "enabled websecurity"

# Row 8 - qtbrowser-2.0.10
$ws.Range("A8").Value = "UPC"
$ws.Range("B8").Value = "DAWN Europe- Browser"
$ws.Range("C8").Value = "qtbrowser-2.0.10"
$ws.Range("D8").Value = "10c09a0589129af0f9b5963f719ed6fa2790e922"
$ws.Range("E8").Value = "Conditionally added syslog functionality"

# Row 9 - qtbrowser-2.0.11
$ws.Range("A9").Value = "UPC"
$ws.Range("B9").Value = "DAWN Europe- Browser"
$ws.Range("C9").Value = "qtbrowser-2.0.11"
$ws.Range("D9").Value = "7193f983d5a07207a67dd22ec0ef0877b885e4fc"
$ws.Range("E9").Value = "Fixed syslog issue"

# New rows inherit the same "vertical-top" style used on A6/B6 for columns A & B
$ws.Range("A7:B9").VerticalAlignment = -4160

# --- Column widths (character units -> XML width = ColumnWidth + 5/6) ---
$ws.Columns.Item(1).ColumnWidth = 10
$ws.Columns.Item(2).ColumnWidth = 20.4962962962963
$ws.Columns.Item(3).ColumnWidth = 12.9518518518519
$ws.Columns.Item(4).ColumnWidth = 40.4074074074074
$ws.Columns.Item(5).ColumnWidth = 98.5333333333333

# --- Hyperlink: update display text to match the full remark text ---
$ws.Range("E6").Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("E6"), "http://download.qt.io/official_releases/qt/5.4/5.4.1/submodules/", [Type]::Missing, [Type]::Missing, "First qt5.4.1 delivery to Pace as part of the S1b build within the Dawn project`n(http://download.qt.io/official_releases/qt/5.4/5.4.1/submodules/)") | Out-Null
# Restore the original (non-hyperlink) formatting that Hyperlinks.Add overwrote
$ws.Range("E6").WrapText = $true
$ws.Range("E6").Font.Name = "Calibri"
$ws.Range("E6").Font.Size = 12
$ws.Range("E6").Font.Color = 0
$ws.Range("E6").Font.Underline = $false

# --- Row heights (set last: hyperlink/font ops above auto-fit the row, so
#     the explicit height must be (re)applied after them) ---
$ws.Rows.Item(6).RowHeight = 15.8

# --- Selection / active cell ---
$ws.Range("C10").Select()
